$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.555558562278748
$ws.Range("B1").Value = 2.460288524627686
$ws.Range("C1").Value = 2.827755451202393
$ws.Range("D1").Value = 3.284664869308472
$ws.Range("E1").Value = 1.664823055267334
